$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "CB" (nom) column, shifting
# nom/url_produit one column to the right (CB->CC, CC->CD).
$ws.Columns("CB:CB").Insert()

# Header for the newly inserted column: the next scrape timestamp.
$ws.Range("CB1").Value = "2026-01-31 08:15:32"

# The new column duplicates the last price snapshot (old CB, still in CA)
# for every data row.
$lastRow = 206
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 80).Value = $ws.Cells.Item($r, 79).Value()
}
